$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.858.56'
$ws.Range("E2").Value = '  +1.03%  '
$ws.Range("D3").Value = '3.245.29'
$ws.Range("E3").Value = '  +1.79%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.90'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.84%  '
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").Value = '3.242.39'
$ws.Range("E8").Value = '  +1.85%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.549'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.62%  '
$ws.Range("E10").Value = '  +0.92%  '
$ws.Range("E11").Value = '  -6.75%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.513'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000272'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '39.14'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.67%  '
$ws.Range("D15").Value = '3.776.62'
$ws.Range("E15").Value = '  +1.60%  '
$ws.Range("D16").Value = '66.797.67'
$ws.Range("E16").Value = '  +0.92%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.53'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.61%  '
$ws.Range("D18").Value = '3.241.65'
$ws.Range("E18").Value = '  +1.49%  '
$ws.Range("E19").Value = '  +1.25%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '513.38'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.41'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.81%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.738'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.37%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.12'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.62%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.94'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.39%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.96'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.21%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.35'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.90%  '
$ws.Range("E28").Value = '  +0.86%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.43'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.04'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.51%  '
$ws.Range("E31").Value = '  +1.95%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.39'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.70%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("E34").Value = '  -3.34%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.56'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.96%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '523.04'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.80%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0957'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '56.25'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.65%  '
$ws.Range("D39").Value = '0.0₃0767'
$ws.Range("E39").Value = '  +19.11%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0422'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.98%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.02'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.54%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.129'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.85'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.40%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.303'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.87%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.50'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.09%  '
$ws.Range("D46").Value = '2.870.37'
$ws.Range("E46").Value = '  -1.62%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.63'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.00%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.43'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.58%  '
$ws.Range("E50").Value = '  +0.86%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.62'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.59%  '
